$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Build out the analysis (COUNTIF / percentage) columns on the
#    existing "Confidence Raw" sheet - mirrors the pattern already used
#    on "Post Skill level Raw".
# ------------------------------------------------------------------
$raw = $wb.Worksheets.Item("Confidence Raw")

$raw.Range("E1").Value = "Percent"

$raw.Range("C2").Value = "Strongly agree"
$raw.Range("D2").Formula = "=COUNTIF(A2:A138, ""Strongly agree"" )"
$raw.Range("E2").Formula = "=D2/D7*100"

$raw.Range("C3").Value = "Somewhat agree"
$raw.Range("D3").Formula = "=COUNTIF(A2:A138, ""Somewhat agree"")"
$raw.Range("E3").Formula = "=D3/D7*100"

$raw.Range("C4").Value = "Neither agree nor disagree"
$raw.Range("D4").Formula = "=COUNTIF(A2:A138, ""Neither agree nor disagree"")"
$raw.Range("E4").Formula = "=D4/D7*100"

$raw.Range("C5").Value = "Somewhat disagree"
$raw.Range("D5").Formula = "=COUNTIF(A2:A138, ""Somewhat disagree"")"
$raw.Range("E5").Formula = "=D5/D7*100"

$raw.Range("C6").Value = "Strongly disagree"
$raw.Range("D6").Formula = "=COUNTIF(A2:A138, ""Strongly disagree"")"
$raw.Range("E6").Formula = "=D6/D7*100"

$raw.Range("C7").Value = "total"
$raw.Range("D7").Formula = "=SUM(D2:D6)"

# Selection on the raw sheet reverts to a plain single-sheet view once
# it is no longer the active tab.
$raw.Range("A1").Select()

# ------------------------------------------------------------------
# 2. Insert the new "Confidence Likert" summary sheet right after
#    "Confidence Raw" (mirrors "Post Skill Level Likert" etc.).
# ------------------------------------------------------------------
$likert = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $raw)
$likert.Name = "Confidence Likert"

$likert.Range("A1").Value = "I felt confident when completing today's camp activities."
$likert.Range("A1").Interior.ColorIndex = 15

$likert.Range("C7").Value = "Strongly agree"
$likert.Range("D7").Value = 78
$likert.Range("E7").Value = 56.934306569343065

$likert.Range("C8").Value = "Somewhat agree"
$likert.Range("D8").Value = 39
$likert.Range("E8").Value = 28.467153284671532

$likert.Range("C9").Value = "Neither agree nor disagree"
$likert.Range("D9").Value = 11
$likert.Range("E9").Value = 8.0291970802919703

$likert.Range("C10").Value = "Somewhat disagree"
$likert.Range("D10").Value = 4
$likert.Range("E10").Value = 2.9197080291970803

$likert.Range("C11").Value = "Strongly disagree"
$likert.Range("D11").Value = 5
$likert.Range("E11").Value = 3.6496350364963499

$likert.Range("C12").Value = "total"
$likert.Range("D12").Value = 137

$likert.Range("C7").Select()
$likert.Activate()
